$wb = $excel.ActiveWorkbook

# GlobalConstantFloatTable is the second sheet
$ws = $wb.Worksheets.Item("GlobalConstantFloatTable")

# Add the two new rows
$ws.Range("A15").Value = "TutorialStartX"
$ws.Range("C15").Value = -1.5

$ws.Range("A16").Value = "TutorialStartZ"
$ws.Range("C16").Value = -3

# Select A16 on this sheet and make it the active sheet/tab
$ws.Range("A16").Select()
$ws.Activate()
